$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Derslik" (classroom) table had its "ID" column removed — GA no longer
# takes an explicit ID parameter, it now keys off İsim/Tür/Kapasite only.
# Unlist the table first so the underlying worksheet column delete is a
# plain grid operation (data + shared strings + dimension + column widths
# all shift left by one), then rebuild the table over the new A1:C11 range
# so it keeps covering İsim/Tür/Kapasite with the same style.
$tbl = $ws.ListObjects.Item("Tablo1")
$tblName = $tbl.Name
$tbl.Unlist()

# Remove the whole "ID" column (column A): shifts İsim/Tür/Kapasite left
# into A/B/C and drops the ID values + header entirely.
$ws.Columns.Item(1).Delete()

$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:C11"), $false, $true)
$newTbl.Name = $tblName
$newTbl.TableStyle = "TableStyleLight8"

# Leave the selection on the (now first) column, matching a "select column,
# delete" workflow.
[void]$ws.Columns.Item(1).Select()
